$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.301.35'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '1.921.77'
$ws.Range("E3").Value = '  -0.41%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.38%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7431'
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '244.27'
$ws.Range("E6").Value = '  -2.16%  '
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '27.34'
$ws.Range("E8").Value = '  -2.06%  '
$ws.Range("E9").Value = '  -2.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06985'
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7744'
$ws.Range("E11").Value = '  -1.90%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07995'
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("D13").Value = '1.940.97'
$ws.Range("E13").Value = '  +0.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.309'
$ws.Range("E14").Value = '  -1.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.84'
$ws.Range("E15").Value = '  -2.83%  '
$ws.Range("D16").Value = '30.320.87'
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.26'
$ws.Range("E17").Value = '  -1.54%  '
$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '246.63'
$ws.Range("E18").Value = '  -2.57%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.867'
$ws.Range("E19").Value = '  +1.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007859'
$ws.Range("E20").Value = '  -2.52%  '
$ws.Range("D21").Value = '2.180.81'
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.003'
$ws.Range("E22").Value = '  +0.33%  '
$ws.Range("E23").Value = '  +0.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.672'
$ws.Range("E24").Value = '  -2.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.432'
$ws.Range("E26").Value = '  +0.93%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.00'
$ws.Range("E27").Value = '  -0.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1274'
$ws.Range("E28").Value = '  -4.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.142'
$ws.Range("E29").Value = '  -6.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.362'
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.550'
$ws.Range("E31").Value = '  +1.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.358'
$ws.Range("E32").Value = '  -1.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.078'
$ws.Range("E33").Value = '  -1.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05194'
$ws.Range("E34").Value = '  +1.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.308'
$ws.Range("E35").Value = '  +0.86%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7513'
$ws.Range("E36").Value = '  +0.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.778'
$ws.Range("E37").Value = '  +0.32%  '
$ws.Range("E38").Value = '  -1.57%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.790'
$ws.Range("E39").Value = '  -0.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.421'
$ws.Range("E40").Value = '  +0.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '76.04'
$ws.Range("E41").Value = '  -2.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4467'
$ws.Range("E42").Value = '  -0.84%  '
$ws.Range("E43").Value = '  -1.79%  '
$ws.Range("E44").Value = '  +0.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8383'
$ws.Range("E45").Value = '  -0.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.691'
$ws.Range("E46").Value = '  +2.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.42'
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.860'
$ws.Range("E48").Value = '  +0.81%  '
$ws.Range("D49").Value = '2.077.25'
$ws.Range("E49").Value = '  -0.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.21'
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1224'
$ws.Range("E51").Value = '  +5.50%  '
